$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Settings")
$ws.Activate()

# New exception row (previously blank row 22)
$ws.Range("A22").Value = "TemplateNotFoundException"
$ws.Range("B22").Value = "Template cloud not be selected based on given data"

# Offer letter placeholder keys - renamed (drop the underscore) and a new
# Zip placeholder inserted, Job Title placeholder text tweaked
$ws.Range("A31").Value = "OfferLetterDate"
$ws.Range("B31").Value = "__Date__"

$ws.Range("A32").Value = "OfferLetterName"
$ws.Range("B32").Value = "__Name__"

$ws.Range("A33").Value = "OfferLetterLastName"
$ws.Range("B33").Value = "__LastName__"

$ws.Range("A34").Value = "OfferLetterAddress"
$ws.Range("B34").Value = "__Address__"

$ws.Range("A35").Value = "OfferLetterCity"
$ws.Range("B35").Value = "__City__"

$ws.Range("A36").Value = "OfferLetterState"
$ws.Range("B36").Value = "__State__"

$ws.Range("A37").Value = "OfferLetterZip"
$ws.Range("B37").Value = "__Zip__"

$ws.Range("A38").Value = "OfferLetterJobTitle"
$ws.Range("B38").Value = "__Job Title__"

$ws.Range("A39").Value = "OfferLetterSalary"
$ws.Range("B39").Value = "__Salary__"

$ws.Range("A40").Value = "OfferLetterOutputFolder"
$ws.Range("B40").Value = "C:\Users\55649C\Documents\UiPath\P004_SP002_090_NewHireCommunication_OfferLetterCreation_Performer\Data\Output\OfferLetters\"

$ws.Range("A41").Value = "OfferLetterTemplateFolder"
$ws.Range("B41").Value = "C:\Users\55649C\Documents\UiPath\P004_SP002_090_NewHireCommunication_OfferLetterCreation_Performer\Data\Input\Templates\"

$ws.Range("A42").Value = "Template_Exempt_With_Trial"
$ws.Range("B42").Value = "Exempt - with Trial Period.docx"

$ws.Range("A43").Value = "Template_Exempt_Without_Trial"
$ws.Range("B43").Value = "Exempt - without Trial Period.docx"

$ws.Range("A44").Value = "Template_NonExempt_With_Trial"
$ws.Range("B44").Value = "Nonexempt - with Trial Period.docx"

$ws.Range("A45").Value = "Template_NonExempt_Without_Trial"
$ws.Range("B45").Value = "Nonexempt - without Trial Period.docx"

$ws.Range("A46").Value = "Template_Executive"
$ws.Range("B46").Value = "X 5%.docx"

$ws.Range("A47").Value = "Template_Temporary"
$ws.Range("B47").Value = "X 999 Temp.docx"

# New rows for the HR approval mail + signed-document naming
$ws.Range("A49").Value = "HRApprovalMailSubject"
$ws.Range("B49").Value = "Offer Letter Approval - "

$ws.Range("A50").Value = "OfferLetterSignedExtension"
$ws.Range("B50").Value = "_Signed.docx"

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Activate()
$wsAssets.Columns(2).ColumnWidth = 53.5
$wsAssets.Range("A2").Select()

# ---------------------------------------------------------------------------
# Restore Settings as the active sheet/selection, matching the saved view
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B53").Select()
